$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3128
$ws.Range("I74").Value = 2859.2
$ws.Range("K74").Value = 2859.2
$ws.Range("M74").Value = -1923.2
$ws.Range("H77").Value = 3128
$ws.Range("I77").Value = 2859.2
$ws.Range("K77").Value = 14296
$ws.Range("M77").Value = -9616
$ws.Range("H114").Value = 29999
$ws.Range("J114").Value = 29999
$ws.Range("L114").Value = 29999
$ws.Range("N114").Value = -38677
$ws.Range("H132").Value = 7582073.5
$ws.Range("I132").Value = 11500298
$ws.Range("K132").Value = 34500894
$ws.Range("M132").Value = -34498364

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 860
$ws.Range("I2").Value = 680
$ws.Range("J2").Value = 1040
$ws.Range("K2").Value = 680
$ws.Range("L2").Value = 1040
$ws.Range("M2").Value = -567
$ws.Range("N2").Value = -1266
$ws.Range("H32").Value = 6130.9795
$ws.Range("I32").Value = 4963.727
$ws.Range("K32").Value = 4963.727
$ws.Range("M32").Value = -4676.727
$ws.Range("H61").Value = 52632492
$ws.Range("I61").Value = 62500704
$ws.Range("K61").Value = 62500704
$ws.Range("M61").Value = -62500492
$ws.Range("H116").Value = 860
$ws.Range("I116").Value = 680
$ws.Range("J116").Value = 1040
$ws.Range("K116").Value = 680
$ws.Range("L116").Value = 1040
$ws.Range("M116").Value = 1614
$ws.Range("N116").Value = -5628
$ws.Range("H136").Value = 52632492
$ws.Range("I136").Value = 62500704
$ws.Range("K136").Value = 187502112
$ws.Range("M136").Value = -187499562

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 860
$ws.Range("I3").Value = 680
$ws.Range("J3").Value = 1040
$ws.Range("K3").Value = 680
$ws.Range("L3").Value = 1040
$ws.Range("M3").Value = -566
$ws.Range("N3").Value = -1268
$ws.Range("H80").Value = 511.83334
$ws.Range("J80").Value = 595.86664
$ws.Range("L80").Value = 595.86664
$ws.Range("N80").Value = -2591.86664
$ws.Range("H83").Value = 511.83334
$ws.Range("J83").Value = 595.86664
$ws.Range("L83").Value = 2979.3332
$ws.Range("N83").Value = -12963.3332

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 512.5
$ws.Range("I105").Value = 416.66666
$ws.Range("K105").Value = 416.66666
$ws.Range("M105").Value = 1330.33334
$ws.Range("H107").Value = 851
$ws.Range("I107").Value = 465.0909
$ws.Range("K107").Value = 465.0909
$ws.Range("M107").Value = 1454.9091
$ws.Range("H122").Value = 765.1818
$ws.Range("I122").Value = 832.5
$ws.Range("J122").Value = 585.6667
$ws.Range("K122").Value = 2497.5
$ws.Range("L122").Value = 1757.0001
$ws.Range("M122").Value = -47.5
$ws.Range("N122").Value = -6657.0001
$ws.Range("H132").Value = 1834.5
$ws.Range("I132").Value = 1503.4615
$ws.Range("J132").Value = 2551.75
$ws.Range("K132").Value = 4510.3845
$ws.Range("L132").Value = 7655.25
$ws.Range("M132").Value = -1980.3845
$ws.Range("N132").Value = -12715.25
$ws.Range("H134").Value = 18520140
$ws.Range("I134").Value = 1598.9131
$ws.Range("J134").Value = 125001750
$ws.Range("K134").Value = 4796.7393
$ws.Range("L134").Value = 375005250
$ws.Range("M134").Value = -2261.7393
$ws.Range("N134").Value = -375010320

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1275.5
$ws.Range("H93").Value = 6617.8184
$ws.Range("J93").Value = 6617.8184
$ws.Range("L93").Value = 19853.4552
$ws.Range("N93").Value = -23597.4552
$ws.Range("H96").Value = 7584.615
$ws.Range("I96").Value = 800
$ws.Range("J96").Value = 8150
$ws.Range("K96").Value = 2400
$ws.Range("L96").Value = 24450
$ws.Range("M96").Value = -341
$ws.Range("N96").Value = -28568
$ws.Range("H115").Value = 4693.6665
$ws.Range("I115").Value = 2721
$ws.Range("K115").Value = 8163
$ws.Range("M115").Value = -6988

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 13581
$ws.Range("I62").Value = 12077
$ws.Range("K62").Value = 12077
$ws.Range("M62").Value = -11391
$ws.Range("H65").Value = 13581
$ws.Range("I65").Value = 12077
$ws.Range("K65").Value = 36231
$ws.Range("M65").Value = -32799
$ws.Range("H132").Value = 2587.4814
$ws.Range("I132").Value = 2330.9524
$ws.Range("J132").Value = 3485.3333
$ws.Range("K132").Value = 6992.8572
$ws.Range("L132").Value = 10455.9999
$ws.Range("M132").Value = -4462.8572
$ws.Range("N132").Value = -15515.9999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2256.4348
$ws.Range("I132").Value = 1370
$ws.Range("J132").Value = 2938.3076
$ws.Range("K132").Value = 4110
$ws.Range("L132").Value = 8814.9228
$ws.Range("M132").Value = -1580
$ws.Range("N132").Value = -13874.9228

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 3000
$ws.Range("J31").Value = 3000
$ws.Range("L31").Value = 3000
$ws.Range("N31").Value = -3696
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()
$ws.Range("H42").Value = 3000
$ws.Range("J42").Value = 3000
$ws.Range("L42").Value = 3000
$ws.Range("N42").Value = -3756
$ws.Range("H49").Value = 15000
$ws.Range("J49").Value = 15000
$ws.Range("L49").Value = 15000
$ws.Range("N49").Value = -15460
$ws.Range("H51").Value = 5000
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -6020
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H107").Value = 782
$ws.Range("I107").Value = 825
$ws.Range("J107").Value = 753.3333
$ws.Range("K107").Value = 2475
$ws.Range("L107").Value = 2259.9999
$ws.Range("M107").Value = -555
$ws.Range("N107").Value = -6099.9999
$ws.Range("H113").Value = 346.625
$ws.Range("J113").Value = 471.375
$ws.Range("L113").Value = 1414.125
$ws.Range("N113").Value = -5754.125
$ws.Range("H122").Value = 25001614
$ws.Range("I122").Value = 25001614
$ws.Range("K122").Value = 75004842
$ws.Range("M122").Value = -75002392
$ws.Range("H132").Value = 3359
$ws.Range("I132").Value = 3338.9285
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 10016.7855
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -7486.7855
$ws.Range("N132").Value = -15558.5
$ws.Range("H136").Value = 1207.238
$ws.Range("I136").Value = 1207
$ws.Range("J136").Value = 1207.7142
$ws.Range("K136").Value = 3621
$ws.Range("L136").Value = 3623.1426
$ws.Range("M136").Value = -1071
$ws.Range("N136").Value = -8723.142599999999
$ws.Range("H139").Value = 34280
$ws.Range("J139").Value = 34280
$ws.Range("L139").Value = 34280
$ws.Range("N139").Value = -44560
